$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 452 (weekly Acelga price data at
# Mercado Mayorista Lo Valledor de Santiago gets appended, pushing the
# previously-last 10 rows down by 3).
$ws.Rows("452:454").Insert()

# Shared/static values for every detail row in this sheet.
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$unidad    = "`$/docena de atados"
$origen    = "Región Metropolitana"
$clasif    = "Hortaliza"

$newRows = @(
    @{ Row = 452; Calidad = "Extra";   Vol = 170; Pmin = 12000; Pmax = 12000; Pprom = 12000; PKg = 4000 },
    @{ Row = 453; Calidad = "Primera"; Vol = 200; Pmin = 10000; Pmax = 10000; Pprom = 10000; PKg = 3333 },
    @{ Row = 454; Calidad = "Segunda"; Vol = 130; Pmin = 8000;  Pmax = 8000;  Pprom = 8000;  PKg = 2667 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 6
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 44448
    $ws.Cells.Item($row, 5).Value  = 13
    $ws.Cells.Item($row, 6).Value  = 100112009
    $ws.Cells.Item($row, 7).Value  = "Acelga"
    $ws.Cells.Item($row, 8).Value  = "Sin especificar"
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.Pmin
    $ws.Cells.Item($row, 12).Value = $r.Pmax
    $ws.Cells.Item($row, 13).Value = $r.Pprom
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = 3
    $ws.Cells.Item($row, 18).Value = $clasif
}
